$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Feuil1")

# Existing agent "Léa" (row 3): zone reassigned, EPI/Extincteur flags flipped
$ws.Range("B3").Value = "PRM11"
$ws.Range("F3").Value = "NON"
$ws.Range("G3").Value = "OUI"

# Expand the table to welcome a new agent row
$tbl = $ws.ListObjects.Item("Tableau1")
$tbl.Resize($ws.Range("A1:G5"))

# New agent "Brann"
$ws.Range("A4").Value = "Brann"
$ws.Range("B4").Value = "Brann"
$ws.Range("C4").Value = "04h - 12h "
$ws.Range("D4").Value = "A"
$ws.Range("E4").Value = "OUI"
$ws.Range("F4").Value = "NON"
$ws.Range("G4").Value = "OUI"

# Underline the Extincteur flag for the new row
$ws.Range("G4").Font.Underline = $true

$ws.Range("B4").Select()
